$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.160.53'
$ws.Cells.Item(2, 5).Value = '  -0.59%  '

$ws.Cells.Item(3, 4).Value = '1.825.16'

$ws.Cells.Item(4, 4).Value = '''0.9996'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.08%  '

$ws.Cells.Item(5, 4).Value = '''234.77'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.80%  '

$ws.Cells.Item(6, 4).Value = '''0.5994'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -4.19%  '

$ws.Cells.Item(7, 5).Value = '  +0.04%  '

$ws.Cells.Item(8, 5).Value = '  -4.88%  '

$ws.Cells.Item(9, 4).Value = '''0.2789'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -3.59%  '

$ws.Cells.Item(10, 4).Value = '''23.43'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -5.40%  '

$ws.Cells.Item(12, 4).Value = '1.828.42'
$ws.Cells.Item(12, 5).Value = '  -0.55%  '

$ws.Cells.Item(13, 4).Value = '''4.788'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -3.27%  '

$ws.Cells.Item(14, 4).Value = '''0.6284'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -6.73%  '

$ws.Cells.Item(15, 4).Value = '''0.000009887'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -3.06%  '

$ws.Cells.Item(16, 4).Value = '2.079.11'
$ws.Cells.Item(16, 5).Value = '  -0.41%  '

$ws.Cells.Item(17, 4).Value = '''78.82'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -3.47%  '

$ws.Cells.Item(18, 4).Value = '''5.846'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -5.81%  '

$ws.Cells.Item(19, 4).Value = '29.180.99'
$ws.Cells.Item(19, 5).Value = '  -0.58%  '

$ws.Cells.Item(20, 4).Value = '''226.02'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.62%  '

$ws.Cells.Item(21, 5).Value = '  +0.01%  '

$ws.Cells.Item(22, 4).Value = '''11.70'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -4.72%  '

$ws.Cells.Item(23, 4).Value = '''6.981'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -4.89%  '

$ws.Cells.Item(24, 5).Value = '  +0.01%  '

$ws.Cells.Item(25, 4).Value = '''155.13'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.88%  '

$ws.Cells.Item(26, 4).Value = '''8.013'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -5.32%  '

$ws.Cells.Item(27, 4).Value = '''0.1297'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -3.40%  '

$ws.Cells.Item(28, 4).Value = '''16.54'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -4.52%  '

$ws.Cells.Item(29, 4).Value = '''1.483'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.85%  '

$ws.Cells.Item(30, 4).Value = '''0.06238'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -14.28%  '

$ws.Cells.Item(31, 4).Value = '''1.448'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.93%  '

$ws.Cells.Item(32, 5).Value = '  -5.14%  '

$ws.Cells.Item(33, 5).Value = '  -5.87%  '

$ws.Cells.Item(34, 4).Value = '''1.120'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.64%  '

$ws.Cells.Item(35, 4).Value = '''1.738'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -4.33%  '

$ws.Cells.Item(36, 4).Value = '''0.6396'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -7.99%  '

$ws.Cells.Item(37, 4).Value = '''2.534'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.32%  '

$ws.Cells.Item(38, 4).Value = '1.215.89'
$ws.Cells.Item(38, 5).Value = '  -1.15%  '

$ws.Cells.Item(39, 4).Value = '''2.729'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.94%  '

$ws.Cells.Item(40, 4).Value = '''0.01730'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -5.51%  '

$ws.Cells.Item(41, 4).Value = '''6.496'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -6.24%  '

$ws.Cells.Item(42, 4).Value = '''0.9046'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -4.31%  '

$ws.Cells.Item(43, 4).Value = '''1.000'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.02%  '

$ws.Cells.Item(44, 4).Value = '1.986.47'
$ws.Cells.Item(44, 5).Value = '  -0.20%  '

$ws.Cells.Item(45, 4).Value = '''100.24'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -0.29%  '

$ws.Cells.Item(46, 5).Value = '  -4.05%  '

$ws.Cells.Item(47, 4).Value = '''0.00000000117'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.13%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).Value = '''8.519'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -4.21%  '

$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).Value = '''1.595'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -6.64%  '

$ws.Cells.Item(50, 4).Value = '''0.4551'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.70%  '

$ws.Cells.Item(51, 4).Value = '''0.05500'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -2.65%  '
